# Apply BOM updates to Airbrake Motor Controller workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Description text updates (shared strings) ---
# Row 20: "Diode 10TQ035" -> "D Schottky", "Schottky Rectifier" -> "Schottky Diode"
$ws.Range("A20").Value = "D Schottky"
$ws.Range("B20").Value = "Schottky Diode"

# --- Updated supplier unit price / subtotal values ---
$ws.Range("G6").Value = 0.0257
$ws.Range("H6").Value = 0.3855

$ws.Range("G15").Value = 0.0621
$ws.Range("H15").Value = 0.0621

$ws.Range("G18").Value = 0.0105
$ws.Range("H18").Value = 0.0105

$ws.Range("G20").Value = 0.4242
$ws.Range("H20").Value = 0.8484

$ws.Range("G26").Value = 0.1346
$ws.Range("H26").Value = 0.1346

$ws.Range("G29").Value = 3.08
$ws.Range("H29").Value = 3.08

$ws.Range("G34").Value = 0.0331
$ws.Range("H34").Value = 0.0331
